# Updates LearnKneeLocalization.xlsx: rename the worksheet tab to match
# the workbook file name, refresh the Element Name / Description table
# (rows 2-18) with the corrected / reshuffled anatomy entries (old sounds
# -related placeholder rows are gone, new knee-anatomy rows were added),
# fix two rows that had accidentally inherited the bold header style, and
# move the saved selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was still named after the Polish Excel default ("Arkusz1");
# rename it to match the workbook file name.
$ws.Name = "LearnKneeLocalization"

$ws.Cells.Item(2, 1).Value = "Rzepka kolanowa"
$ws.Cells.Item(2, 2).Value = "Jest to kość spłaszczona trójkątna z zaokrąglonymi brzegami, włączona w ścięgno mięśnia czworogłowego uda i położona od dolnego końca kości udowej."

$ws.Cells.Item(3, 1).Value = "Kość piszczelowa"
$ws.Cells.Item(3, 2).Value = "Jest kością o funkcji strukturalnej, tworzącą podudzie. Należy do kości długich. Koniec dalszy kości piszczelowej tworzy kostkę przyśrodkową. "

$ws.Cells.Item(4, 1).Value = "Kość udowa"
$ws.Cells.Item(4, 2).Value = "To najdłuższa kość ciała ludzkiego. Jak każda kość długa składa się z trzonu i dwóch końców. Na powierzchni tylnej trzonu zaznacza się wydatna kresa chropowa (linea aspera), złożona z dwóch warg, przyśrodkowej (labium mediale) i bocznej (labium laterale). Ku dołowi obie wargi rozchodząc się ograniczają powierzchnię podkolanową (facies poplitea). Warga przyśrodkowa kończy się u dołu guzkiem przywodziciela (tuberculum adductorium). "

$ws.Cells.Item(5, 1).Value = "Więzadło rzepki"
$ws.Cells.Item(5, 2).Value = "Więzadło rzepki łączy ją z kością piszczelową w ruchach zgięcia i prostowania stawu kolanowego rzepka porusza się razem z kością piszczelową. W przypadku wyprostowanego kolana i napiętych mięśni prostujących część dolna powierzchni stawowej rzepki spoczywa na powierzchni rzepkowej kości udowej, część górna leży powyżej; gdy kolano jest zgięte, rzepka przesuwa się ku dołowi i do tyłu, układa się w rowku między obu kłykciami kości udowej i jest unieruchomiona. Gdy kolano jest wyprostowane i mięśnie rozkurczone, rzepka może być przesuwana na boki. Poniżej powierzchni stawowej znajduje się wypukłe chropowate pole, którego dolna część służy za przyczep więzadła rzepki."
$ws.Cells.Item(5, 1).Font.Bold = $false

$ws.Cells.Item(6, 1).Value = "Więzadło krzyżowe przednie"
$ws.Cells.Item(6, 2).Value = "Przebiega od powierzchni przyśrodkowej kłykcia bocznego kości udowej przednio-przyśrodkowo do pola międzykłykciowego przedniego kości piszczelowej"

$ws.Cells.Item(7, 1).Value = "Więzadło krzyżowe tylne"
$ws.Cells.Item(7, 2).Value = "Przebiega od powierzchni przyśrodkowej kłykcia przyśrodkowego do dołu, do pola międzykłykciowego tylnego kości piszczelowej"

$ws.Cells.Item(8, 1).Value = "Więzadło poprzeczne"
$ws.Cells.Item(8, 2).Value = "Więzadło poprzeczne kolana rozpięte jest między najbardziej do przodu położonymi punktami obu łąkotek, łącząc je. Jest to więzadło cienkie, okrągławe, często powstrzymane w rozwoju. Nieraz może go całkowicie brakować. Napina się ono podczas ruchów obrotowych stawu na zewnątrz."

$ws.Cells.Item(9, 1).Value = "Łąkotka przyśrodkowa"
$ws.Cells.Item(9, 2).Value = "Półkolista, w kształcie litery C.  Przyczepia się do pola międzykłykciowego przedniego i tylnego. Mocno przytwierdzona przyśrodkowo do torebki stawowej, a bocznie - do więzadła pobocznego piszczelowego, które ogranicza jej ruchomość"
$ws.Rows.Item(9).RowHeight = 57.6

$ws.Cells.Item(10, 1).Value = "Łąkotka boczna"
$ws.Cells.Item(10, 2).Value = "Prawie całkowicie kolista. Przyczpia się do pola międzykłykciowego przedniego i tylnego. Nie jest przytwierdzona do torebki stawowej, więc jest bardziej ruchoma od łąkotki przyśrodkowej"
$ws.Rows.Item(10).RowHeight = 43.2

$ws.Cells.Item(11, 1).Value = "Kłykieć przyśrodkowy"
$ws.Cells.Item(11, 2).Value = "Kłykieć przyśrodkowy stawu kolanowego, zwany również guzkiem przyśrodkowym lub wyrostkiem przyśrodkowym, to ważna struktura anatomiczna w obrębie stawu kolanowego. Jest to wyrostek kostny znajdujący się na wewnętrznej stronie kości udowej, blisko stawu kolanowego. Kłykieć przyśrodkowy jest kluczowy dla stabilności i funkcji stawu kolanowego "
$ws.Rows.Item(11).RowHeight = 86.4

$ws.Cells.Item(12, 1).Value = "Kość strzałkowa"
$ws.Cells.Item(12, 2).Value = "Znajduje się po stronie przyśrodkowej goleni i po kości udowej jest najdłuższą kością szkieletu. W przekroju poprzecznym ma kształt trójkątny. U góry, gdzie bierze udział w wytwarzaniu stawu kolanowgo, kość znacznie grubieje, ku dołowi zwęża się, a następnie znowu poszerza, choć w mniejszym stopniu niż u góry. Jak każda kość długa składa się z trzonu i dwóch końców"
$ws.Rows.Item(12).RowHeight = 86.4

$ws.Cells.Item(13, 1).Value = "Powierzchnia stawowa piszczela"
$ws.Cells.Item(13, 2).Value = "Powierzchnia stawowa kości piszczelowej jest wklęsła i przylega do odpowiedniej powierzchni stawowej kości udowej. Kształt tej powierzchni jest asymetryczny, co pozwala na pewne ograniczenia w ruchomości stawu kolanowego, chroniąc go przed nadmiernymi skręceniami."
$ws.Rows.Item(13).RowHeight = 72

$ws.Cells.Item(14, 1).Value = "Powierzchnia stawowa rzepki"
$ws.Cells.Item(14, 2).Value = "Czyli tylna, powyżej wierzchołka pokryta jest grubą warstwą chrząstki szklistej dzieli się ona na dwa pola, z których pole boczne jest większe od przyśrodkowego obie te części przedzielone są podłużnie biegnącym wzniesieniem. Odpowiada ono podłużnemu rowkowi na powierzchni rzepkowej kości udowej."
$ws.Rows.Item(14).RowHeight = 72

$ws.Cells.Item(15, 1).Value = "Powierzchnia stawowa strzałki"
$ws.Cells.Item(15, 2).Value = "Płaska powierzchnia stawowa strzałkowa zlokalizowana jest na kłykciu bocznym kości piszczelowej. Powierzchnia stawowa strzałkowa skierowana jest ku dołowi, ku tyłowi i bocznie i przylega do powierzchni stawowej głowy strzałki."
$ws.Rows.Item(15).RowHeight = 57.6
$ws.Cells.Item(15, 1).Font.Bold = $false

$ws.Cells.Item(16, 1).Value = "Więzadło poboczne strzałkowe"
$ws.Cells.Item(16, 2).Value = "Rozpościera się od bocznej powierzchni kości udowej aż do tzw. głowy kości strzałkowej i odpowiada za stabilność kolana od strony bocznej"
$ws.Rows.Item(16).RowHeight = 43.2

$ws.Cells.Item(17, 1).Value = "Więzadło poboczne piszczelowe"
$ws.Cells.Item(17, 2).Value = "Rozpościera się od przyśrodkowej powierzchni kości udowej aż do przyśrodkowej powierzchni piszczeli, stąd odpowiada za stabilność stawu kolanowego od strony przyśrodkowej (czyli od wewnątrz)."

$ws.Cells.Item(18, 1).Value = "Więzadło łąkotkowo-udowe tylne"
$ws.Cells.Item(18, 2).Value = "Więzadło przyczepia się w okolicy tylnego przyczepu łąkotki bocznej. Biegnie ku górze i przyśrodkowo do tyłu od więzadła krzyżowego tylnego. Przeważnie łączy się z nim kończąc się w miejscu jego przyczepu do wewnętrznej powierzchni kłykcia przyśrodkowego kości udowej."
$ws.Rows.Item(18).RowHeight = 72

# The saved cursor/selection moves from A16 to D2, and the view no longer
# needs to remember a scrolled-down top-left cell (previously A10).
$ws.Range("D2").Select()

